$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PayNowData")

# Row 6 is a copy of Row 2, except column A (Notes) and column B (ID)
$ws.Range("A2:AB2").Copy()
$ws.Range("A6").PasteSpecial(-4104)

$ws.Range("A6").Value = "No Modify Amount"
$ws.Range("B6").Value = "5"

$ws.Rows.Item(6).RowHeight = $ws.Rows.Item(2).RowHeight

$ws.Range("R2").Copy()
$ws.Range("R6").PasteSpecial(-4122)

$ws.Range("C6").Select()
